# Updated cryptos list on Tue Oct  8 20:24:04 UTC 2024 with GitHub Actions
#
# Refreshes the price / 1h-volume-change snapshot for the crypto tracker
# sheet, plus re-ranks row 51 (Stellar -> BabyDogeCoin).
#
# Note: several "Price" (column D) values are plain decimal numbers
# rendered as text (e.g. "19.50", "0.531") so that trailing zeros are
# preserved exactly as scraped. Excel's Range.Value setter auto-converts
# any numeric-looking string to a real number, which would both change
# the cell's type and silently drop significant trailing zeros (e.g.
# "19.50" -> 19.5). To keep those cells as text we prefix them with a
# leading apostrophe, Excel's standard "force text" convention, before
# assigning.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.240.53'
$ws.Range('E2').Value = '  -1.48%  '
$ws.Range('D3').Value = '2.442.74'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''583.46'
$ws.Range('E5').Value = '  +2.06%  '
$ws.Range('D6').Value = '''143.93'
$ws.Range('E6').Value = '  -1.89%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '''0.531'
$ws.Range('E8').Value = '  -1.19%  '
$ws.Range('D9').Value = '2.440.86'
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  -3.34%  '
$ws.Range('E11').Value = '  +2.60%  '
$ws.Range('E12').Value = '  -1.09%  '
$ws.Range('D13').Value = '''0.345'
$ws.Range('E13').Value = '  -3.17%  '
$ws.Range('D14').Value = '''26.46'
$ws.Range('E14').Value = '  -1.97%  '
$ws.Range('E15').Value = '  -3.59%  '
$ws.Range('D16').Value = '2.864.23'
$ws.Range('E16').Value = '  -0.99%  '
$ws.Range('D17').Value = '62.111.60'
$ws.Range('E17').Value = '  -1.61%  '
$ws.Range('D18').Value = '2.430.57'
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('D19').Value = '''10.90'
$ws.Range('E19').Value = '  -3.55%  '
$ws.Range('D20').Value = '''7.13'
$ws.Range('E20').Value = '  -2.69%  '
$ws.Range('D21').Value = '''330.26'
$ws.Range('E21').Value = '  +0.65%  '
$ws.Range('E22').Value = '  -2.51%  '
$ws.Range('E23').Value = '  -4.43%  '
$ws.Range('E24').Value = '  -3.92%  '
$ws.Range('D25').Value = '''65.70'
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('D26').Value = '''9.38'
$ws.Range('E26').Value = '  +4.81%  '
$ws.Range('D27').Value = '''619.58'
$ws.Range('E27').Value = '  +0.79%  '
$ws.Range('D28').Value = '2.566.68'
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('D29').Value = '0.0₃0957'
$ws.Range('E29').Value = '  -7.05%  '
$ws.Range('E30').Value = '  -0.14%  '
$ws.Range('E31').Value = '  -4.47%  '
$ws.Range('E32').Value = '  -3.03%  '
$ws.Range('E33').Value = '  -0.17%  '
$ws.Range('E34').Value = '  -0.51%  '
$ws.Range('E35').Value = '  -5.06%  '
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('E37').Value = '  -6.27%  '
$ws.Range('D38').Value = '''0.379'
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('D39').Value = '''151.58'
$ws.Range('E39').Value = '  +2.72%  '
$ws.Range('E40').Value = '  -2.36%  '
$ws.Range('E41').Value = '  -3.04%  '
$ws.Range('E42').Value = '  -1.35%  '
$ws.Range('D43').Value = '''42.42'
$ws.Range('E43').Value = '  +1.36%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('E45').Value = '  -6.31%  '
$ws.Range('D46').Value = '''143.34'
$ws.Range('E46').Value = '  -3.58%  '
$ws.Range('E47').Value = '  -3.38%  '
$ws.Range('E48').Value = '  -1.12%  '
$ws.Range('E49').Value = '  -0.33%  '
$ws.Range('D50').Value = '''19.50'
$ws.Range('E50').Value = '  -7.86%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0235'
$ws.Range('E51').Value = '  +7.03%  '
